# Eine Thangka raus genommen und ein paar Preise geändert
#
# - Move the "Shakyamuni als Kusumari" row (sheet1 row 52) to the
#   "catalog_not_found" sheet (row 18), marking it as not-available (E=0).
# - Move the "Das Leben von Shakyamuni" thangka row (sheet1 row 167) to the
#   "catalog_not_found" sheet (row 19), keeping its availability flag.
# - Bump a handful of prices on catalog_product_bearbeitet.
# - Leave the UI selection/active-sheet state matching the author's last
#   on-screen position.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # catalog_product_bearbeitet
$ws2 = $wb.Worksheets.Item(2)   # catalog_not_found

# --- Remove the Buddha statue "Shakyamuni als Kusumari" (row 52) ---------
# Cut (keeps formatting) to the first free row on the "not found" sheet …
$ws1.Range("A52:G52").Cut($ws2.Range("A18")) | Out-Null
# … then wipe what remains on the source sheet so the row disappears
# entirely instead of leaving an empty, still-styled row behind.
$ws1.Range("A52:G52").Clear() | Out-Null
# It is no longer being sold, so flag it as unavailable.
$ws2.Range("E18").Value = 0

# --- Remove the Thangka "Das Leben von Shakyamuni" (row 167) -------------
$ws1.Range("A167:G167").Cut($ws2.Range("A19")) | Out-Null
$ws1.Range("A167:G167").Clear() | Out-Null

# --- Price updates on catalog_product_bearbeitet --------------------------
$ws1.Range("D160").Value = 450
$ws1.Range("D161").Value = 500
$ws1.Range("D162").Value = 500
$ws1.Range("D166").Value = 500

# --- Restore the on-screen selection / active sheet state -----------------
$ws2.Range("B24").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("D161").Select() | Out-Null
